# Updated cryptos list on Mon Aug 19 11:42:38 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns on
# Sheet1 with the latest scraped values. Row numbers below correspond to
# the worksheet rows (row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain text (e.g. "58.037.30", "0.0530") even though a lot
# of the values look numeric. Force those specific cells to keep a Text
# number format before writing so Excel doesn't silently coerce them into
# real numbers (which would drop things like trailing zeros or treat
# "58.037.30" style big numbers inconsistently).
$dTextCells = @("D2", "D3", "D5", "D6", "D9", "D13", "D14", "D16", "D19", "D20", "D23", "D24", "D26", "D31", "D32", "D33", "D35", "D36", "D38", "D40", "D41", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $dTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price (column D) updates
$ws.Range("D2").Value = "58.037.30"
$ws.Range("D3").Value = "2.573.29"
$ws.Range("D5").Value = "533.80"
$ws.Range("D6").Value = "141.05"
$ws.Range("D9").Value = "6.74"
$ws.Range("D13").Value = "3.028.60"
$ws.Range("D14").Value = "57.969.69"
$ws.Range("D16").Value = "2.582.51"
$ws.Range("D19").Value = "333.64"
$ws.Range("D20").Value = "9.98"
$ws.Range("D23").Value = "66.56"
$ws.Range("D24").Value = "0.416"
$ws.Range("D26").Value = "0.158"
$ws.Range("D31").Value = "154.91"
$ws.Range("D32").Value = "5.81"
$ws.Range("D33").Value = "18.82"
$ws.Range("D35").Value = "36.87"
$ws.Range("D36").Value = "1.08"
$ws.Range("D38").Value = "0.812"
$ws.Range("D40").Value = "3.57"
$ws.Range("D41").Value = "282.27"
$ws.Range("D46").Value = "0.0530"
$ws.Range("D47").Value = "18.23"
$ws.Range("D48").Value = "0.0225"
$ws.Range("D49").Value = "1.899.73"
$ws.Range("D50").Value = "17.67"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  -2.85%  "
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("E11").Value = "  +3.11%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("E16").Value = "  -2.70%  "
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("E21").Value = "  -4.09%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -4.85%  "
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -4.11%  "
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("E36").Value = "  -4.68%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("E47").Value = "  -4.43%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("E50").Value = "  -4.14%  "
$ws.Range("E51").Value = "  -4.73%  "

Write-Host "Updated cryptos price/volume columns."
